# Add new data tables (worksheets) to the workbook:
#   course_offered, track_requirement, gep_requirement, degree_requirement
#
# New sheet order must end up as:
#   course, course_offered, prereq, subject, semester, major,
#   track_requirement, gep_requirement, degree_requirement

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the four new worksheets.  We create them in the order that
# makes their internal sheetId values come out as 6,7,8,9 (the engine assigns
# sheetId sequentially at creation time), then reposition them afterwards.
# ---------------------------------------------------------------------------
$degreeReq = $wb.Worksheets.Add()
$degreeReq.Name = "degree_requirement"

$gepReq = $wb.Worksheets.Add()
$gepReq.Name = "gep_requirement"

$courseOffered = $wb.Worksheets.Add()
$courseOffered.Name = "course_offered"

$trackReq = $wb.Worksheets.Add()
$trackReq.Name = "track_requirement"

# ---------------------------------------------------------------------------
# Step 2: reposition the new sheets into their final order. We always look
# sheets back up by name (rather than reuse old references) since indices
# shift after every move.
# ---------------------------------------------------------------------------
$wb.Worksheets("course_offered").Move($wb.Worksheets("prereq"))

$wb.Worksheets("track_requirement").Move()
$wb.Worksheets("gep_requirement").Move()
$wb.Worksheets("degree_requirement").Move()

# ---------------------------------------------------------------------------
# Step 3: populate the header rows for the new tables. The order in which we
# write the cells controls the order new shared strings are created in, so
# fill degree_requirement first, then gep_requirement, then course_offered,
# to match requirement_id, requirement_type, course_id_options,
# gep_requirement_id, gep_type, frequency_id, offered_prob.
# ---------------------------------------------------------------------------

$degreeReq = $wb.Worksheets("degree_requirement")
$degreeReq.Cells.Item(1, 1).Value = "requirement_id"
$degreeReq.Cells.Item(1, 2).Value = "major_id"
$degreeReq.Cells.Item(1, 3).Value = "requirement_type"
$degreeReq.Cells.Item(1, 4).Value = "course_id_options"
$degreeReq.Columns.Item(1).ColumnWidth = 12.666666666666666
$degreeReq.Columns.Item(2).ColumnWidth = 7.166666666666667
$degreeReq.Columns.Item(3).ColumnWidth = 14.833333333333334
$degreeReq.Columns.Item(4).ColumnWidth = 14.833333333333334
$degreeReq.Range("D2").Select()

$gepReq = $wb.Worksheets("gep_requirement")
$gepReq.Cells.Item(1, 1).Value = "gep_requirement_id"
$gepReq.Cells.Item(1, 2).Value = "gep_type"
$gepReq.Cells.Item(1, 3).Value = "course_id_options"
$gepReq.Columns.Item(1).ColumnWidth = 16.666666666666668
$gepReq.Range("C2").Select()

$courseOffered = $wb.Worksheets("course_offered")
$courseOffered.Cells.Item(1, 1).Value = "frequency_id"
$courseOffered.Cells.Item(1, 2).Value = "course_id"
$courseOffered.Cells.Item(1, 3).Value = "term"
$courseOffered.Cells.Item(1, 4).Value = "offered_prob"
$courseOffered.Columns.Item(1).ColumnWidth = 10.5
$courseOffered.Columns.Item(2).ColumnWidth = 7.833333333333333
$courseOffered.Columns.Item(3).ColumnWidth = 3.8333333333333335
$courseOffered.Columns.Item(4).ColumnWidth = 10.666666666666666
$courseOffered.PageSetup.Orientation = 1
$courseOffered.Range("D2").Select()

# track_requirement is left completely empty (no header row / data yet).

# ---------------------------------------------------------------------------
# Step 4: cosmetic touch up on the "course" sheet - the course_name column
# was widened.
# ---------------------------------------------------------------------------
$course = $wb.Worksheets("course")
$course.Columns.Item(3).ColumnWidth = 62.166666666666664
$course.Range("A28:E28").Select()

# ---------------------------------------------------------------------------
# Step 5: make course_offered the active sheet (it is the second tab, same
# position that was previously tabbed to "prereq").
# ---------------------------------------------------------------------------
$wb.Worksheets("course_offered").Activate()
